$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is a new data row appended after the existing data (rows 2-52).
# The date-like text must stay as plain text (matching the rest of the
# "date" column, which is stored as text, not as a real Excel date), so a
# leading apostrophe is used to suppress Excel's automatic date detection;
# the style is then reset to "Normal" so no stray quote-prefix / text
# number-format style gets attached to the cell.
$ws.Range("A53").Value = "'2025/10/02"
$ws.Range("A53").Style = "Normal"

$ws.Range("B53").Value = "木"
$ws.Range("C53").Value = 20
$ws.Range("D53").Value = 201
